$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.371.81"
$ws.Range("E2").Value = "  +1.20%  "

$ws.Range("D3").Value = "3.010.58"
$ws.Range("E3").Value = "  +0.18%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "507.53"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -0.89%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.42"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +0.41%  "

$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("E8").Value = "  -0.18%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.56"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  +0.45%  "

$ws.Range("E10").Value = "  +0.74%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.365"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +2.42%  "

$ws.Range("D12").Value = "3.527.12"
$ws.Range("E12").Value = "  +0.20%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.129"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  +0.08%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.32"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +2.32%  "

$ws.Range("E15").Value = "  +2.72%  "

$ws.Range("D16").Value = "57.408.99"
$ws.Range("E16").Value = "  +1.16%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.20"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  +4.03%  "

$ws.Range("D18").Value = "3.008.67"
$ws.Range("E18").Value = "  +0.07%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.79"
$ws.Range("D19").NumberFormat = "General"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.94"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +1.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "327.54"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -0.98%  "

$ws.Range("E22").Value = "  -0.15%  "

$ws.Range("E23").Value = "  -1.88%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.498"
$ws.Range("D24").NumberFormat = "General"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.33"
$ws.Range("D25").NumberFormat = "General"

$ws.Range("E26").Value = "  -3.57%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -0.06%  "

$ws.Range("D28").Value = "0.0₃0916"
$ws.Range("E28").Value = "  +1.32%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.74"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +0.19%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.32"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +3.78%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.80"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -0.14%  "

$ws.Range("E32").Value = "  -5.00%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.53"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -0.79%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.75"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +4.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "153.88"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +0.18%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.87"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +3.46%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.27"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -0.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.64"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +3.69%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0675"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -0.34%  "

$ws.Range("D40").Value = "3.043.81"
$ws.Range("E40").Value = "  +0.22%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "37.80"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +2.33%  "

$ws.Range("E42").Value = "  +4.74%  "

$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.649"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -0.10%  "

$ws.Range("E45").Value = "  -0.45%  "

$ws.Range("D46").Value = "2.222.60"
$ws.Range("E46").Value = "  -2.44%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.977"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -2.99%  "

$ws.Range("E48").Value = "  +3.66%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0238"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -1.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.47"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +0.22%  "

$ws.Range("E51").Value = "  -5.16%  "
